$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 38

# Row 3
$ws.Range("C3").Value = 44

# Row 4
$ws.Range("B4").Value = "<hot>"
$ws.Range("C4").Value = 46

# Row 5
$ws.Range("B5").Value = "<this>"
$ws.Range("C5").Value = 35

# Row 6
$ws.Range("C6").Value = 34

# Row 7
$ws.Range("C7").Value = 39

# Row 8
$ws.Range("B8").Value = "<number>"
$ws.Range("C8").Value = 43

# Row 9
$ws.Range("B9").Value = "<hotel>"
$ws.Range("C9").Value = 43

# Row 10
$ws.Range("B10").Value = "<sae>"
$ws.Range("C10").Value = 37

# Row 11
$ws.Range("B11").Value = "<copa>"
$ws.Range("C11").Value = 36

# Row 12
$ws.Range("C12").Value = 37

# Row 13
$ws.Range("C13").Value = 41

# Row 14
$ws.Range("B14").Value = "<nomem>"
$ws.Range("C14").Value = 38

# Row 15
$ws.Range("C15").Value = 39

# Row 16
$ws.Range("C16").Value = 40

# Row 17
$ws.Range("B17").Value = "<escape>"
$ws.Range("C17").Value = 36

# Row 18
$ws.Range("B18").Value = "<with>"
$ws.Range("C18").Value = 31
